$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$enc1 = 'U2FsdGVkX18GwQhBbrS0mtak6vzynMD8Hu2/fDkarVoyt0+Tr/IOTIj8oLOOZn74Ldb7kQzWhEFw40b1VWY8B8rF7tZmyiwksHf/ul/SVbyrDZF4enLa3LQ1bML/Zp1W56aKIwp3R2m08zE8OO5WNdofBVm/ELy+HCi+ZfuE11iI/J0PxAa5gmjSdI2gpwS3g8qs+/bliQGE32dOMkEHRgYZqhZzZe3OO/lH2ilKvvqj9mGLpTMpvwCKVfH3hsUdU8JeBj3+8HNF9SilKTrxK1icrcXAX5/UPznNmGHEo1WHhBb0VsUt2OTKKxYf656pNT+eRfTEQKA/0JbGac/hH36a8tguW0nzPCBNKp6+gGZNN2NQtG/0NTXr2vJHC3UPBpunN5fUGMPM8EDFH452lVpXhkCBiU7xdC+7cEWZ7B5usqnfSFVtHJVoNReC2n0J8N9ljuEIB+d5z6tSmBS33RWOHLZEf+h125R4KYWuMQm5jCh4YNupKZX83K60mMCgwsIqz9BAzuRC2iJ+wBaF7NrEJO0p8YDnO3vyI6sWinplGVKgVkQKdo+rcoJspJsVNSCE6BasuYU42KWVMkUzkUKeEs/FWmC+/u9hZoqmfNgxZ4VWfWkezhdb3sqo1RZ7jyrFzTFlW9SPuPU9J4ueESYfvah2a0dbvuw+5jLs2PdayFkrosdMqUPs6KTDnax+J76XTtqTanbJ+d5FKaX2ZAVeEgxUP/7F6KoYFbbEebgkOGplx8fXkmXSMuQYoSEYiqKVE7BVI8qyQA7QB/TY8Q2wx+7Ns/2Ngi7pSyKQQDavZ6Ra4OnedpETxxR6W8Kb8dS9/wuTi6ymx8t8tVkogZFqOttoVlJwb3aNRTr8frg47WHTzmMH2fiSXwP6N9KUiyuiOJ0LRxEQVXhCue84cgXmWp/cuUfFmuqBV2amc3nyNoMdv4G2JmSzJiBQpE3ASnxsVnAJz+S3rfY/OjV4k1Bg+i33ZDfSsaWXejjip/6vVECRyrXS6neVfWTMKpgWGcjgKw+5ba6vKNSbtk8O/k5EYgLS9F+yC2gyZlx9jUSwNX2024kmPMOaaiaee6cVGujJuAhclqluGbwxqL2AFbo78O91kE+2ip8JofbTQsg2OWE4M7iwpty0acjalA1/xJpQgipQQfsQIi2ou4paI4H/JzvkV3GnRJECSu7zUPYVy5kYH16Mjdc/fXZ5nj/gie92116zm1BGcceTGr7id2h/7cXsf25woHKqt6ketg3UjxPBpzfHMt4Jzr5mZ7IaQaV3t3efddZv14md9FH2fyZeeOCw+ehxbepQNzBO9fo='
$enc2 = 'U2FsdGVkX18CeviQ24hIEERMC5lllAkCyTq8qB+K1J8ELLtppxyKsUFxOE+nntq9/WW38YXJQgmprt6+xWNAmDhQSQC6lExckmwooXEjTNnJ3+9TvopXVeVD4S0/efSymJkoILmDy9RhiqXoSqODsFf1mw71a6OLivpxueg/q+qqYdvq9yTI5iNQnz8y5ZvghwX8Tnm229QpdNMizeUnveiSMqY7iXdMtRUQvaeCOZiQmAcsaONQgiZHeeIJfpeKgFkOVDkheWibIj5j8MfxlLI/fFNKAEvF/SJnCDktvpPdSOUz0PHB2+E7GdBSsjK5'

$ws.Range("B2").Value2 = $enc1
$ws.Range("C2").Value2 = "2:55-3:0"

$ws.Range("B3").Value2 = $enc2
$ws.Range("C3").Value2 = "3:0-3:5"
$ws.Range("C3").HorizontalAlignment = -4108

$ws.Range("B4").Value2 = $enc1
$ws.Range("C4").Value2 = "10:55-11:0"
$ws.Range("C4").HorizontalAlignment = -4108

$ws.Range("B5").Value2 = $enc2
$ws.Range("C5").Value2 = "11:0-11:5"
$ws.Range("C5").HorizontalAlignment = -4108

$ws.Range("B6").Value2 = $enc1
$ws.Range("C6").Value2 = "18:55-19:0"
$ws.Range("C6").HorizontalAlignment = -4108

$ws.Range("B7").Value2 = $enc2
$ws.Range("C7").Value2 = "19:0-19:5"
$ws.Range("C7").HorizontalAlignment = -4108

$ws.Range("B12").Select()
